$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 14 (A14) used to be a link formula `=[1]tot_tut!A20` showing "ToT-TuT".
# That underlying external row was renamed/repurposed, so A14 now just
# holds the literal label "ASG" (computation of lenders' profit / ASG row).
$ws.Range("A14").Value = 'ASG'

# Rows 24/25 hypothesis-test labels: ToT-TuT -> ASG
$ws.Range("A24").Value = "`$H_0`$ : ASG=0"
$ws.Range("A25").Value = "`$H_0`$ : ASG`$\geq`$ 0"

# Reset the view: no frozen/scrolled top-left cell, selection anchored at A2
# (top of the A2:E25 table) instead of E25.
$ws.Range("A2:E25").Select() | Out-Null
